$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 20 de Septiembre de 2020 a las 18:41'

$data = @(
    @('Estados Unidos',6983461,16058,4235275,2544220,0,142,203966),
    @('India',5418681,20451,4314606,1017149,0,152,86926),
    @('Brasil',4531539,3192,3820095,574767,0,112,136677),
    @('Rusia',1103399,6148,909357,174624,0,79,19418),
    @('Peru',762865,0,607837,123659,0,0,31369),
    @('Colombia',758398,0,627685,106674,0,0,24039),
    @('Mexico',694121,5167,496224,124639,0,455,73258),
    @('Sudafrica',659656,0,589434,54282,0,0,15940),
    @('España',659334,0,0,0,0,0,30495),
    @('Argentina',622934,0,488231,121794,0,110,12909),
    @('Chile',444674,0,418101,14319,0,0,12254),
    @('Francia',442194,0,91574,319346,0,0,31274),
    @('Iran',422140,3097,359570,38269,0,183,24301),
    @('Reino Unido',394257,3899,0,0,0,18,41777),
    @('Banglades',348918,1544,256565,87414,0,26,4939),
    @('Arabia Saudita',329754,483,310439,14830,0,27,4485),
    @('Irak',319035,3438,253591,56889,0,64,8555),
    @('Pakistan',305671,640,292303,6952,0,1,6416),
    @('Turquia',302867,1519,266117,29244,0,61,7506),
    @('Italia',298156,1587,218351,44098,0,15,35707),
    @('Filipinas',286743,3311,229865,51894,0,55,4984),
    @('Alemania',272670,362,243500,19704,0,0,9466),
    @('Indonesia',244676,3989,177327,57796,0,105,9553),
    @('Israel',183602,0,132449,49927,0,0,1226),
    @('Ucrania',175678,2966,77512,94609,0,41,3557),
    @('Canada',143603,829,124684,9702,0,6,9217),
    @('Bolivia',130470,419,89032,33852,0,36,7586),
    @('Ecuador',125620,0,97063,17473,0,0,11084),
    @('Catar',123376,230,120303,2863,0,1,210),
    @('Rumania',112781,1231,89771,18575,0,33,4435),
    @('Republica Dominicana',107700,0,80820,24836,0,0,2044),
    @('Kazajistan',107262,63,101877,3714,0,0,1671),
    @('Panama',105601,0,80190,23164,0,0,2247),
    @('Egipto',101900,0,88666,7484,0,0,5750),
    @('Belgica',100748,1099,18945,71859,0,7,9944),
    @('Marruecos',99816,0,79008,19013,0,0,1795),
    @('Kuwait',99434,385,90168,8682,0,3,584),
    @('Paises Bajos',93778,1844,0,0,0,4,6279),
    @('Oman',93475,1722,85418,7211,0,28,846),
    @('Suecia',88237,0,0,0,0,0,5865),
    @('China',85279,10,80477,168,0,0,4634),
    @('Guatemala',85152,0,74497,7550,0,0,3105),
    @('Emiratos Arabes Unidos',84916,674,74273,10239,0,0,404),
    @('Polonia',79240,910,64302,12645,0,11,2293),
    @('Japon',78073,0,70495,6083,0,0,1495),
    @('Bielorrusia',75461,0,73212,1473,0,0,776),
    @('Honduras',71143,532,21810,47167,0,20,2166),
    @('Portugal',68577,552,45596,21069,0,13,1912),
    @('Etiopia',68131,0,27939,39103,0,0,1089),
    @('Venezuela',65949,0,55155,10255,0,0,539),
    @('Barein',64499,0,57299,6979,0,0,221),
    @('Nepal',64122,1325,46233,17478,0,10,411),
    @('Costa Rica',63712,0,23552,39454,0,0,706),
    @('Singapur',57576,18,57181,368,0,0,27),
    @('Nigeria',57145,0,48431,7619,0,0,1095),
    @('Uzbekistan',51235,243,47271,3535,0,2,429),
    @('Argelia',49623,0,34923,13035,0,0,1665),
    @('Suiza',49283,0,40500,6738,0,0,2045),
    @('Chequia',48881,575,24744,23635,0,3,502),
    @('Armenia',47431,277,42608,3893,0,2,930),
    @('Moldavia',46596,260,34236,11157,0,2,1203),
    @('Ghana',46004,127,45153,554,0,0,297),
    @('Kirguistan',45335,0,41484,2788,0,0,1063),
    @('Azerbaiyan',39188,146,36755,1858,0,1,575),
    @('Afganistan',39044,125,32576,5027,0,4,1441),
    @('Austria',38095,621,29229,8100,0,1,766),
    @('Kenia',36981,152,23887,12446,0,2,648),
    @('Estado de Palestina',35686,683,23700,11724,0,9,262),
    @('Paraguay',33015,0,17535,14844,0,0,636),
    @('Irlanda',32933,395,23364,7777,0,0,1792),
    @('Serbia',32908,68,31411,756,0,1,741),
    @('Libano',29303,1006,12047,16959,0,11,297),
    @('Libia',27949,715,15068,12437,0,8,444),
    @('El Salvador',27553,125,21561,5181,0,3,811),
    @('Australia',26898,13,24062,1987,0,5,849),
    @('Bosnia y Herzegovina',25428,211,17878,6787,0,5,763),
    @('Corea del Sur',22975,82,20158,2434,0,5,383),
    @('Dinamarca',22905,469,17514,4753,0,3,638),
    @('Camerun',20431,0,19124,891,0,0,416),
    @('Costa de Marfil',19269,0,18392,757,0,0,120),
    @('Bulgaria',18819,0,13558,4506,0,0,755),
    @('Hungria',17990,1070,4391,12916,0,8,683),
    @('Republica de Macedonia',16735,178,13863,2179,0,4,693),
    @('Madagascar',16053,33,14646,1184,0,4,223),
    @('Grecia',15142,164,9989,4815,0,7,338),
    @('Croacia',14922,197,12536,2138,0,4,248),
    @('Senegal',14714,26,11260,3152,0,0,302),
    @('Zambia',14131,61,13365,436,0,0,330),
    @('Sudan',13535,0,6759,5940,0,0,836),
    @('Noruega',12883,25,10371,2245,0,0,267),
    @('Albania',12385,159,6940,5083,0,4,362),
    @('Consejo Danes para los Refugiados',10515,27,9930,314,0,3,271),
    @('Namibia',10377,85,8033,2232,0,1,112),
    @('Guinea',10286,0,9681,542,0,0,63),
    @('Malasia',10219,52,9355,734,0,0,130),
    @('Tunez',9736,626,2386,7195,0,17,155),
    @('Guayana Francesa',9692,0,9341,286,0,0,65),
    @('Maldivas',9649,0,8188,1428,0,0,33),
    @('Tayikistan',9346,43,8108,1165,0,0,73),
    @('Finlandia',8980,58,7700,941,0,0,339),
    @('Gabon',8696,0,7848,795,0,0,53),
    @('Haiti',8600,0,6363,2016,0,0,221),
    @('Montenegro',7898,0,5129,2635,0,0,134),
    @('Luxemburgo',7804,0,6703,977,0,0,124),
    @('Zimbabue',7672,0,5914,1533,0,0,225),
    @('Mauritania',7365,0,6927,277,0,0,161),
    @('Mozambique',6771,234,3622,3106,0,2,43),
    @('Eslovaquia',6677,131,3548,3090,0,0,39),
    @('Uganda',6287,270,2616,3608,0,0,63),
    @('Malaui',5718,0,4030,1509,0,0,179),
    @('Birmania',5541,671,1260,4189,0,11,92),
    @('Republica de Yibuti',5403,0,5333,9,0,0,61),
    @('Suazilandia',5245,0,4571,570,0,0,104),
    @('Cabo Verde',5186,0,4581,555,0,0,50),
    @('Cuba',5091,36,4412,564,0,2,115),
    @('Hong Kong',5033,23,4708,222,0,0,103),
    @('Guinea Ecuatorial',5002,0,4509,410,0,0,83),
    @('Congo',4986,0,3887,1010,0,0,89),
    @('Nicaragua',4961,0,2913,1901,0,0,147),
    @('Republica de Africa Central',4786,0,1830,2894,0,0,62),
    @('Jordania',4779,239,2844,1905,0,0,30),
    @('Jamaica',4758,0,1327,3371,0,0,60),
    @('Surinam',4709,0,4383,229,0,0,97),
    @('Ruanda',4689,0,2910,1753,0,0,26),
    @('Eslovenia',4420,111,3023,1255,0,1,142),
    @('Angola',3901,0,1445,2309,0,0,147),
    @('Trinidad yTobago',3866,13,1753,2049,0,3,64),
    @('Siria',3765,0,932,2663,0,0,170),
    @('Lituania',3744,80,2198,1459,0,0,87),
    @('Mayotte',3541,0,2964,537,0,0,40),
    @('Gambia',3526,22,1992,1426,0,0,108),
    @('Tailandia',3506,6,3340,107,0,0,59),
    @('Georgia',3502,196,1494,1989,0,0,19),
    @('Aruba',3460,0,2128,1309,0,0,23),
    @('Somalia',3442,41,2877,467,0,0,98),
    @('Guadalupe',3426,0,837,2563,0,0,26),
    @('Sri Lanka',3287,4,3088,186,0,0,13),
    @('Bahamas',3214,0,1679,1462,0,0,73),
    @('Reunion',3194,0,2482,697,0,0,15),
    @('Mali',3006,0,2349,529,0,0,128),
    @('Estonia',2924,49,2377,483,0,0,64),
    @('Malta',2731,32,2047,664,0,1,20),
    @('Sudan del Sur',2642,0,1290,1303,0,0,49),
    @('Botsuana',2567,0,624,1930,0,0,13),
    @('Islandia',2346,39,2121,215,0,0,10),
    @('Guinea-Bisau',2303,0,1127,1137,0,0,39),
    @('Benin',2280,0,1950,290,0,0,40),
    @('Guyana',2168,0,1331,773,0,0,64),
    @('Sierra Leona',2159,0,1650,437,0,0,72),
    @('Yemen',2026,0,1221,220,0,0,585),
    @('Uruguay',1904,0,1612,246,0,0,46),
    @('Burkina Faso',1816,0,1176,584,0,0,56),
    @('Nueva Zelanda',1815,4,1719,71,0,0,25),
    @('Togo',1659,0,1259,359,0,0,41),
    @('Belice',1606,16,876,710,0,0,20),
    @('Republica de Chipre',1590,0,1282,286,0,0,22),
    @('Principado de Andorra',1564,0,1164,347,0,0,53),
    @('Letonia',1525,10,1248,241,0,0,36),
    @('Lesoto',1390,0,754,603,0,0,33),
    @('Liberia',1335,0,1216,37,0,0,82),
    @('Polinesia Francesa',1271,0,1028,241,0,0,2),
    @('Niger',1183,0,1104,10,0,0,69),
    @('Republica del Chad',1149,0,966,102,0,0,81),
    @('Martinica',1122,0,98,1006,0,0,18),
    @('Vietnam',1068,0,942,91,0,0,35),
    @('Santo Tome y Principe',908,0,876,17,0,0,15),
    @('San Marino',723,0,669,12,0,0,42),
    @('Crucero',712,0,651,48,0,0,13),
    @('Islas Turcas y Caicos',668,1,572,91,0,0,5),
    @('San Martin (Parte Holandesa)',584,10,488,76,0,0,20),
    @('Papua Nueva Guinea',516,0,232,278,0,0,6),
    @('Tanzania',509,0,183,305,0,0,21),
    @('Taiwan',507,1,479,21,0,0,7),
    @('Burundi',473,0,374,98,0,0,1),
    @('Comoras',470,0,450,13,0,0,7),
    @('Islas Feroe',431,0,412,19,0,0,0),
    @('Mauricio',366,0,338,18,0,0,10),
    @('Eritrea',364,0,305,59,0,0,0),
    @('Gibraltar',350,0,323,27,0,0,0),
    @('Isla de Man',339,0,312,3,0,0,24),
    @('San Martin (Parte Francesa)',330,0,206,118,0,0,6),
    @('Mongolia',312,1,302,10,0,0,0),
    @('Camboya',275,0,274,1,0,0,0),
    @('Butan',259,1,190,69,0,0,0),
    @('Curazao',247,0,88,158,0,0,1),
    @('Islas Caimanes',208,0,204,3,0,0,1),
    @('Monaco',192,0,152,39,0,0,1),
    @('Barbados',185,0,172,6,0,0,7),
    @('Bermudas',180,0,166,5,0,0,9),
    @('Brunei',145,0,142,0,0,0,3),
    @('Seychelles',141,0,136,5,0,0,0),
    @('Liechtenstein',113,0,109,3,0,0,1),
    @('Antigua y Barbuda',96,0,92,1,0,0,3),
    @('Islas Virgenes Britanicas',69,0,48,20,0,0,1),
    @('San Vicente y las Granadinas',64,0,64,0,0,0,0),
    @('Macao',46,0,46,0,0,0,0),
    @('Puerto Rico',39,0,1,36,0,0,2),
    @('Bonaire, San Eustaquio y Saba',36,0,17,18,0,0,1),
    @('Guam',32,0,0,31,0,0,1),
    @('Fiyi',32,0,26,4,0,0,2),
    @('Timor Oriental',27,0,26,1,0,0,0),
    @('Santa Lucia',27,0,26,1,0,0,0),
    @('Nueva Caledonia',26,0,26,0,0,0,0),
    @('Dominica',24,0,18,6,0,0,0),
    @('Granada',24,0,24,0,0,0,0),
    @('San Bartolome',23,0,16,7,0,0,0),
    @('Laos',23,0,22,1,0,0,0),
    @('Islas Virgenes de los Estados Unidos',17,0,0,17,0,0,0),
    @('San Cristobal y Nieves',17,0,17,0,0,0,0),
    @('Groenlandia',14,0,14,0,0,0,0),
    @('Islas Malvinas',13,0,13,0,0,0,0),
    @('Montserrat',13,0,12,0,0,0,1),
    @('Santa Sede',12,0,12,0,0,0,0),
    @('San Pedro y Miquelon',11,0,5,6,0,0,0),
    @('Sahara Occidental',10,0,8,1,0,0,1),
    @('Anguila',3,0,3,0,0,0,0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

Write-Output "Done, last row = $($r - 1)"
